$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.227.94'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.65%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.353.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.61'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.84%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.62'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.350.87'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.37%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.468'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.16%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.46'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.53%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.95%  '

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.917.94'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.47%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.43%  '

$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.86'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.91%  '

$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000172'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.96%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.349.84'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.286.18'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.00'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.44%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.38%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.23'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.47'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.73%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.551'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.82%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.485.36'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.38%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.07'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000123'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.30%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.78'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +10.16%  '

$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.40%  '

$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.165'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.13'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.23%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.13'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.27%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.50'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.21'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.39%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.75'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.25%  '

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.52'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.87%  '

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.72'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0757'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.84%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.23%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.768'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.60%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.70'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.33%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.40'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.20'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.32%  '

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.38'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.81'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.00%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.79'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.98'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.28%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.349.71'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.20%  '
